# This script applies the diff described:
# - For each data row (2-19, 21-34; row 20 already correct), the D,E,F,G,H
#   (open/close/high/low price and shares_outstanding) values are overwritten
#   with the STNE (row 20) values, and the I column (fixed_ticker) is set to
#   "STNE" for every row, replacing the various other-ticker values.
# - This also causes the now-unused shared strings (CYBR, EBAY, FTNT, ... NTAP)
#   to no longer be referenced; the engine will compact/garbage-collect the
#   shared string table automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "2" = @{ D = 32; E = 28.69000053405762; F = 32.5; G = 27.59000015258789; H = 250599325 }
    "3" = @{ D = 32; E = 28.69000053405762; F = 32.5; G = 27.59000015258789; H = 250599325 }
    "4" = @{ D = 32; E = 28.69000053405762; F = 32.5; G = 27.59000015258789; H = 250599325 }
    "5" = @{ D = 32; E = 28.69000053405762; F = 32.5; G = 27.59000015258789; H = 250599325 }
    "6" = @{ D = 32; E = 28.69000053405762; F = 32.5; G = 27.59000015258789; H = 250599325 }
    "7" = @{ D = 32; E = 28.69000053405762; F = 32.5; G = 27.59000015258789; H = 250599325 }
    "8" = @{ D = 18.29000091552734; E = 22.1200008392334; F = 23.40999984741211; G = 17.98999977111816; H = 250599325 }
    "9" = @{ D = 41; E = 28.79000091552734; F = 45.61999893188477; G = 24.51000022888184; H = 250599325 }
    "10" = @{ D = 30; E = 35.0099983215332; F = 37.56999969482422; G = 28.04000091552734; H = 250599325 }
    "11" = @{ D = 34.72999954223633; E = 36.79000091552734; F = 36.88999938964844; G = 30.93000030517578; H = 250599325 }
    "12" = @{ D = 40.7400016784668; E = 43.43999862670898; F = 45.47999954223633; G = 39.90999984741211; H = 250599325 }
    "13" = @{ D = 20.90999984741211; E = 26.3799991607666; F = 28.53000068664551; G = 17.71999931335449; H = 250599325 }
    "14" = @{ D = 38.93999862670898; E = 47.70999908447266; F = 47.86999893188477; G = 36.70999908447266; H = 250599325 }
    "15" = @{ D = 53.79999923706055; E = 52.54000091552734; F = 60.7599983215332; G = 50.40000152587891; H = 250599325 }
    "16" = @{ D = 85; E = 71.90000152587891; F = 86.36000061035156; G = 68.45999908447266; H = 250599325 }
    "17" = @{ D = 63.18000030517578; E = 64.63999938964844; F = 70.56999969482422; G = 61.70999908447266; H = 250599325 }
    "18" = @{ D = 67.27999877929688; E = 58.84000015258789; F = 70.73999786376953; G = 53.18999862670898; H = 250599325 }
    "19" = @{ D = 34.95000076293945; E = 33.86000061035156; F = 40.04000091552734; G = 30.70999908447266; H = 250599325 }
    "21" = @{ D = 11.81999969482422; E = 9.420000076293944; F = 12.82999992370606; G = 9.060000419616699; H = 250599325 }
    "22" = @{ D = 7.650000095367432; E = 9.579999923706056; F = 9.909999847412109; G = 7.199999809265137; H = 250599325 }
    "23" = @{ D = 9.810000419616699; E = 10.5; F = 12.25; G = 9.609999656677246; H = 250599325 }
    "24" = @{ D = 9.470000267028809; E = 11.15999984741211; F = 11.77000045776367; G = 8.420000076293945; H = 250599325 }
    "25" = @{ D = 9.5; E = 12.31999969482422; F = 12.52999973297119; G = 8.729999542236328; H = 250599325 }
    "26" = @{ D = 12.77000045776367; E = 14.48999977111816; F = 14.6899995803833; G = 11.07999992370606; H = 250599325 }
    "27" = @{ D = 10.61999988555908; E = 9.909999847412109; F = 10.82999992370606; G = 9.340000152587891; H = 250599325 }
    "28" = @{ D = 17.76000022888184; E = 17.19000053405762; F = 18.8700008392334; G = 16.14999961853027; H = 250599325 }
    "29" = @{ D = 16.68000030517578; E = 15.60000038146973; F = 17.93000030517578; G = 15.11999988555908; H = 250599325 }
    "30" = @{ D = 12.02000045776367; E = 13.11999988555908; F = 13.78999996185303; G = 11.61999988555908; H = 250599325 }
    "31" = @{ D = 11.1899995803833; E = 11.10000038146973; F = 11.80000019073486; G = 10.64000034332275; H = 250599325 }
    "32" = @{ D = 8.069999694824219; E = 9.170000076293944; F = 9.689999580383301; G = 7.71999979019165; H = 250599325 }
    "33" = @{ D = 10.47999954223633; E = 14.0600004196167; F = 14.47999954223633; G = 9.829999923706056; H = 250599325 }
    "34" = @{ D = 16.03000068664551; E = 12.77999973297119; F = 16.68000030517578; G = 12.72000026702881; H = 250599325 }
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item([int]$r, 4).Value = $vals.D
    $ws.Cells.Item([int]$r, 5).Value = $vals.E
    $ws.Cells.Item([int]$r, 6).Value = $vals.F
    $ws.Cells.Item([int]$r, 7).Value = $vals.G
    $ws.Cells.Item([int]$r, 8).Value = $vals.H
    $ws.Cells.Item([int]$r, 9).Value = "STNE"
}
